$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.573.11"
$ws.Range("E2").Value = "  -6.10%  "
$ws.Range("D3").Value = "2.206.65"
$ws.Range("E3").Value = "  -6.90%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.30"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.53"
$ws.Range("E6").Value = "  -10.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  -8.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  -10.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.29"
$ws.Range("E10").Value = "  -12.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.98"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0820"
$ws.Range("E12").Value = "  -10.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.71"
$ws.Range("E13").Value = "  -9.17%  "
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.858"
$ws.Range("E15").Value = "  -12.54%  "
$ws.Range("D16").Value = "2.540.44"
$ws.Range("E16").Value = "  -6.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.03"
$ws.Range("E17").Value = "  -8.07%  "
$ws.Range("D18").Value = "2.190.94"
$ws.Range("E18").Value = "  -7.33%  "
$ws.Range("D19").Value = "42.449.16"
$ws.Range("E19").Value = "  -6.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.58"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "0.0₃0953"
$ws.Range("E21").Value = "  -10.33%  "
$ws.Range("E22").Value = "  -13.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.02"
$ws.Range("E23").Value = "  -11.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.13"
$ws.Range("E24").Value = "  -10.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.72"
$ws.Range("E25").Value = "  -9.53%  "
$ws.Range("E26").Value = "  -8.68%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -10.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.15"
$ws.Range("E30").Value = "  -15.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.32"
$ws.Range("E31").Value = "  -9.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0870"
$ws.Range("E32").Value = "  -10.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.62"
$ws.Range("E33").Value = "  -10.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "154.22"
$ws.Range("E34").Value = "  -8.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.76"
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.16"
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("E37").Value = "  +12.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("E39").Value = "  -8.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  -13.58%  "
$ws.Range("E41").Value = "  -6.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0322"
$ws.Range("E42").Value = "  -9.23%  "
$ws.Range("D43").Value = "1.857.29"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.06"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.33"
$ws.Range("E46").Value = "  -11.48%  "
$ws.Range("E47").Value = "  -11.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.38"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.35"
$ws.Range("E49").Value = "  -7.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.73"
$ws.Range("E50").Value = "  -14.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.58"
$ws.Range("E51").Value = "  -6.64%  "
